# Populate the "AI" sheet (sheet1) and the "医療" sheet (sheet11) of the
# response.xlsx data-storage workbook with the scraped article URLs, and
# restore the window/tab selection state that results from that work
# (the "AI" tab ends up active/selected, with the first empty row A4
# highlighted; the "医療" tab's selection moves down to its first empty
# row, A3, too).

$wb = $excel.ActiveWorkbook

$sheet1Urls = @(
    "https://response.jp/article/2019/12/06/329547.html",
    "https://response.jp/article/2019/07/09/324248.html",
    "https://response.jp/article/2019/11/30/329348.html",
    "https://response.jp/article/2019/11/09/328553.html",
    "https://response.jp/article/2019/04/18/321534.html",
    "https://response.jp/article/2019/06/08/323236.html",
    "https://response.jp/article/2019/09/20/326717.html",
    "https://response.jp/article/2019/10/08/327359.html",
    "https://response.jp/article/2019/06/04/323073.html",
    "https://response.jp/article/2018/02/24/306467.html",
    "https://response.jp/article/2019/03/26/320577.html",
    "https://response.jp/article/2019/04/16/321415.html",
    "https://response.jp/article/2019/03/15/320170.html",
    "https://response.jp/article/2019/12/04/329454.html",
    "https://response.jp/special/recent/3532/%E4%BA%BA%E5%B7%A5%E7%9F%A5%E8%83%BD%EF%BC%88AI%EF%BC%89",
    "https://response.jp/article/2019/08/02/325155.html",
    "https://response.jp/article/2019/07/23/324728.html",
    "https://response.jp/article/2019/12/08/329583.html",
    "https://response.jp/article/2019/10/11/327509.html",
    "https://response.jp/article/2017/06/28/296736.html"
)

$sheet11Urls = @(
    "https://response.jp/article/2019/02/27/319571.html",
    "https://response.jp/article/2018/02/19/306228.html",
    "https://response.jp/article/2019/11/25/329129.html",
    "https://response.jp/article/2019/05/14/322286.html",
    "https://response.jp/article/2019/11/11/328602.html",
    "https://response.jp/article/2019/04/03/320884.html",
    "https://response.jp/article/2018/12/21/317416.html",
    "https://response.jp/article/2018/08/31/313531.html",
    "https://response.jp/article/2014/09/10/232104.html",
    "https://response.jp/article/2014/09/10/232114.html",
    "https://response.jp/article/2015/10/19/262285.html",
    "https://response.jp/release/kyodonews_kaigai/20150428/17526.html",
    "https://response.jp/article/2013/10/21/209017.html",
    "https://response.jp/article/2011/09/22/162695.html",
    "https://response.jp/article/2018/07/25/312329.html",
    "https://response.jp/release/kyodonews_kaigai/20140723/11836.html",
    "https://response.jp/release/kyodonews_kaigai/20191011/60431.html",
    "https://response.jp/release/kyodonews_kaigai/20130710/3307.html",
    "https://response.jp/article/2015/11/04/263565.html"
)

function Fill-Urls($sheet, $urls) {
    $row = 3
    foreach ($u in $urls) {
        $cell = $sheet.Cells.Item($row, 1)
        $cell.Value = $u
        # New rows keep the sheet's default (unstyled) formatting instead
        # of inheriting the column's style index.
        [void]($cell.Style = "Normal")
        $row = $row + 1
    }
}

$wsIryou = $wb.Worksheets.Item("医療")
Fill-Urls $wsIryou $sheet11Urls

$wsAI = $wb.Worksheets.Item("AI")
Fill-Urls $wsAI $sheet1Urls

# "医療" selection settles on its first empty row after the new data.
[void]$wsIryou.Activate()
[void]$wsIryou.Range("A3").Select()

# "AI" becomes the active/selected sheet, with the selection on its own
# first empty row.
[void]$wsAI.Activate()
[void]$wsAI.Range("A4").Select()
